# Entitlement export: add a new "national_id" column.
#
# The test fixture's sheet already has an "extra/unexpected" trailing
# column (M) used to exercise the "unexpected column" import-validation
# test. This change inserts a brand-new column M ("national_id") and
# pushes the old column M (the unexpected-column data) one place to the
# right, to N - without touching any of the other existing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M - this shifts the former M column (and its
# styling/formatting) one position to the right, to N, and widens the
# used range from M to N.
$ws.Columns("M").Insert()

# Header for the freshly inserted column.
$ws.Range("M1").Value = "national_id"

# Sample data: only the second data row (row 3) carries a national id
# in the fixture, matching the original "unexpected column" sample data
# layout (row 2 stays blank for this column, like the rest of that row).
$ws.Range("M3").Value = "ABC123456"
